# Update the existing "data" sheet's F-column timestamps (panel query refreshed)
$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

$dataWs.Range("F2").Value = "2021-10-05 14:20:22.138713"
$dataWs.Range("F3").Value = "2021-10-05 14:20:22.138721"
$dataWs.Range("F4").Value = "2021-10-05 14:20:22.138724"
$dataWs.Range("F5").Value = "2021-10-05 14:20:22.138727"

# Add the new "metadata" sheet, placed after "data"
$metaWs = $wb.Worksheets.Add($null, $dataWs)
$metaWs.Name = "metadata"

# Header row (bold, thin border, centered/top aligned - matches "data" sheet header style)
$metaWs.Range("B1").Value = "data_name"
$metaWs.Range("C1").Value = "data_id"
$metaWs.Range("D1").Value = "data_version"
$metaWs.Range("E1").Value = "data_version_created"
$metaWs.Range("F1").Value = "panel_query_time"
$metaWs.Range("G1").Value = "panel_get_request"

$headerRange = $metaWs.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row
$metaWs.Range("A2").Value = 0
$metaWs.Range("A2").Font.Bold = $true
$metaWs.Range("A2").Borders.LineStyle = 1
$metaWs.Range("A2").Borders.Weight = 2
$metaWs.Range("A2").HorizontalAlignment = -4108
$metaWs.Range("A2").VerticalAlignment = -4160

$metaWs.Range("B2").Value = "Familial tumoral calcinosis"
$metaWs.Range("C2").Value = 552
$metaWs.Range("D2").Value = "'1.7"
$metaWs.Range("E2").Value = "2021-03-02T16:12:47.558174Z"
$metaWs.Range("F2").Value = "2021-10-05 14:20:22.134956"
$metaWs.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/552/?format=json"
